$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: worksheet row number, new DAMSLTag (col I), new DialogAct (col J)
$changes = @(
    @{Row=3; I='ba'; J='Appreciation'}
    @{Row=14; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=22; I='ba'; J='Appreciation'}
    @{Row=43; I='sd'; J='Statement-non-opinion'}
    @{Row=57; I='ba'; J='Appreciation'}
    @{Row=64; I='sv'; J='Statement-opinion'}
    @{Row=96; I='sd'; J='Statement-non-opinion'}
    @{Row=97; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=104; I='ba'; J='Appreciation'}
    @{Row=113; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=115; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=119; I='ba'; J='Appreciation'}
    @{Row=122; I='aa'; J='Agree/Accept'}
    @{Row=126; I='sd'; J='Statement-non-opinion'}
    @{Row=144; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=147; I='sd'; J='Statement-non-opinion'}
    @{Row=148; I='sd'; J='Statement-non-opinion'}
    @{Row=151; I='sd'; J='Statement-non-opinion'}
    @{Row=153; I='sd'; J='Statement-non-opinion'}
    @{Row=160; I='%'; J='Uninterpretable'}
    @{Row=176; I='sd'; J='Statement-non-opinion'}
    @{Row=185; I='aa'; J='Agree/Accept'}
    @{Row=189; I='sd'; J='Statement-non-opinion'}
    @{Row=191; I='sv'; J='Statement-opinion'}
    @{Row=208; I='sv'; J='Statement-opinion'}
    @{Row=211; I='ba'; J='Appreciation'}
    @{Row=218; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=220; I='sd'; J='Statement-non-opinion'}
    @{Row=221; I='sd'; J='Statement-non-opinion'}
    @{Row=224; I='aa'; J='Agree/Accept'}
    @{Row=243; I='aa'; J='Agree/Accept'}
    @{Row=246; I='sd'; J='Statement-non-opinion'}
    @{Row=254; I='sv'; J='Statement-opinion'}
    @{Row=255; I='sv'; J='Statement-opinion'}
    @{Row=263; I='aa'; J='Agree/Accept'}
    @{Row=277; I='aa'; J='Agree/Accept'}
    @{Row=303; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=330; I='sd'; J='Statement-non-opinion'}
    @{Row=332; I='sv'; J='Statement-opinion'}
    @{Row=341; I='sd'; J='Statement-non-opinion'}
    @{Row=355; I='ba'; J='Appreciation'}
    @{Row=356; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=370; I='sv'; J='Statement-opinion'}
    @{Row=400; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=403; I='aa'; J='Agree/Accept'}
    @{Row=404; I='sd'; J='Statement-non-opinion'}
    @{Row=417; I='sd'; J='Statement-non-opinion'}
    @{Row=418; I='sd'; J='Statement-non-opinion'}
    @{Row=422; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=424; I='aa'; J='Agree/Accept'}
    @{Row=425; I='sd'; J='Statement-non-opinion'}
    @{Row=426; I='%'; J='Uninterpretable'}
    @{Row=457; I='sd'; J='Statement-non-opinion'}
    @{Row=469; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=478; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=483; I='aa'; J='Agree/Accept'}
    @{Row=493; I='sd'; J='Statement-non-opinion'}
    @{Row=497; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=501; I='sv'; J='Statement-opinion'}
    @{Row=509; I='sv'; J='Statement-opinion'}
    @{Row=510; I='qy'; J='Yes-No-Question'}
    @{Row=512; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=517; I='sv'; J='Statement-opinion'}
    @{Row=524; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=533; I='ba'; J='Appreciation'}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.I
    $ws.Cells.Item($change.Row, 10).Value = $change.J
}

